$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D->E, old E->F) and mirror C's width/style.
$ws.Columns("D").Insert()
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# Populate the new "天威" (TianWei) company column with its data.
$ws.Range("D2").Value = "天威"
$ws.Range("D3").Value = "1200 / 1500"
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = "繳11個月算12個月"
$ws.Range("D7").Value = "無監視器 / 有4隻監視器"
$ws.Range("D8").Value = "最高理賠200倍"
$ws.Range("D19").Value = "阿秋 0933 185 241"

# Update selection to match the saved state.
$ws.Range("E12").Select() | Out-Null
